# Insert a new data row at row 440 (pushes existing rows 440-470 down to 441-471)
# and populate it with the new "Early Majestic" / "Primera" observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(440).Insert()

$ws.Cells.Item(440, 1).Value2  = 10
$ws.Cells.Item(440, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(440, 3).Value2  = "La Araucanía"
$ws.Cells.Item(440, 4).Value2  = 45265
$ws.Cells.Item(440, 5).Value2  = 9
$ws.Cells.Item(440, 6).Value2  = "Fruta"
$ws.Cells.Item(440, 7).Value2  = 100103
$ws.Cells.Item(440, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(440, 9).Value2  = 100103004
$ws.Cells.Item(440, 10).Value2 = "Durazno"
$ws.Cells.Item(440, 11).Value2 = "Early Majestic"
$ws.Cells.Item(440, 12).Value2 = "Primera"
$ws.Cells.Item(440, 13).Value2 = 100
$ws.Cells.Item(440, 14).Value2 = 20000
$ws.Cells.Item(440, 15).Value2 = 20000
$ws.Cells.Item(440, 16).Value2 = 20000
$ws.Cells.Item(440, 17).Value2 = '$/bandeja 18 kilos granel'
$ws.Cells.Item(440, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(440, 19).Value2 = 1111
$ws.Cells.Item(440, 20).Value2 = 18
